$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 365.9355
$ws.Range("J17").Value = 301.51724
$ws.Range("L17").Value = 904.55172
$ws.Range("N17").Value = -1240.55172
$ws.Range("H40").Value = 5462.5
$ws.Range("J40").Value = 6928.5713
$ws.Range("L40").Value = 6928.5713
$ws.Range("N40").Value = -7278.5713
$ws.Range("H43").Value = 4658.8237
$ws.Range("J43").Value = 4213.3335
$ws.Range("L43").Value = 4213.3335
$ws.Range("N43").Value = -4351.3335
$ws.Range("H92").Value = 463.1875
$ws.Range("I92").Value = 319.18182
$ws.Range("J92").Value = 780
$ws.Range("K92").Value = 319.18182
$ws.Range("L92").Value = 780
$ws.Range("M92").Value = 928.81818
$ws.Range("N92").Value = -3276
$ws.Range("H96").Value = 1387.7858
$ws.Range("I96").Value = 943
$ws.Range("K96").Value = 2829
$ws.Range("M96").Value = -1456
$ws.Range("H132").Value = 15179.211
$ws.Range("I132").Value = 1606.1765
$ws.Range("K132").Value = 4818.529500000001
$ws.Range("M132").Value = -2288.529500000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5048.1924
$ws.Range("I32").Value = 5396.9785
$ws.Range("K32").Value = 5396.9785
$ws.Range("M32").Value = -5109.9785
$ws.Range("H74").Value = 4685.143
$ws.Range("I74").Value = 4216
$ws.Range("J74").Value = 7500
$ws.Range("K74").Value = 4216
$ws.Range("L74").Value = 7500
$ws.Range("M74").Value = -3342
$ws.Range("N74").Value = -9248
$ws.Range("H77").Value = 4685.143
$ws.Range("I77").Value = 4216
$ws.Range("J77").Value = 7500
$ws.Range("K77").Value = 21080
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = -16712
$ws.Range("N77").Value = -46236
$ws.Range("H97").Value = 1757.7812
$ws.Range("I97").Value = 1084.826
$ws.Range("J97").Value = 3477.5557
$ws.Range("K97").Value = 1084.826
$ws.Range("L97").Value = 3477.5557
$ws.Range("M97").Value = -588.826
$ws.Range("N97").Value = -4469.5557
$ws.Range("H102").Value = 9475.5
$ws.Range("I102").Value = 4325.7144
$ws.Range("K102").Value = 4325.7144
$ws.Range("M102").Value = -2703.7144
$ws.Range("H122").Value = 4416.143
$ws.Range("I122").Value = 3536.2222
$ws.Range("K122").Value = 10608.6666
$ws.Range("M122").Value = -8158.6666

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4424.1665
$ws.Range("I94").Value = 1118.6
$ws.Range("K94").Value = 1118.6
$ws.Range("M94").Value = -667.5999999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3872.8462
$ws.Range("I31").Value = 2891.5
$ws.Range("J31").Value = 4714
$ws.Range("K31").Value = 2891.5
$ws.Range("L31").Value = 4714
$ws.Range("M31").Value = -2596.5
$ws.Range("N31").Value = -5304
$ws.Range("H34").Value = 3872.8462
$ws.Range("I34").Value = 2891.5
$ws.Range("J34").Value = 4714
$ws.Range("K34").Value = 2891.5
$ws.Range("L34").Value = 4714
$ws.Range("M34").Value = -2689.5
$ws.Range("N34").Value = -5118
$ws.Range("H122").Value = 4411.6875
$ws.Range("J122").Value = 5208.3335
$ws.Range("L122").Value = 15625.0005
$ws.Range("N122").Value = -20525.0005
$ws.Range("H134").Value = 11006
$ws.Range("I134").Value = 10007.2
$ws.Range("K134").Value = 30021.6
$ws.Range("M134").Value = -27486.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 639862
$ws.Range("J37").Value = 639862
$ws.Range("L37").Value = 1919586
$ws.Range("N37").Value = -1919810
$ws.Range("H95").Value = 4990
$ws.Range("I95").Value = 4990
$ws.Range("K95").Value = 14970
$ws.Range("M95").Value = -12911
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 4375
$ws.Range("J102").Value = 4600
$ws.Range("L102").Value = 13800
$ws.Range("N102").Value = -18668
$ws.Range("H106").Value = 4190.9
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 4190.9
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 12572.7
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -14464.7
$ws.Range("H113").Value = 1457.6666
$ws.Range("I113").Value = 2376.5
$ws.Range("J113").Value = 998.25
$ws.Range("K113").Value = 7129.5
$ws.Range("L113").Value = 2994.75
$ws.Range("M113").Value = -4959.5
$ws.Range("N113").Value = -7334.75
$ws.Range("H140").Value = 626660.3
$ws.Range("I140").Value = 668204.3
$ws.Range("K140").Value = 2004612.9
$ws.Range("M140").Value = -1999432.9

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 9745.182000000001
$ws.Range("I97").Value = 649.5
$ws.Range("J97").Value = 20660
$ws.Range("K97").Value = 649.5
$ws.Range("L97").Value = 20660
$ws.Range("M97").Value = -153.5
$ws.Range("N97").Value = -21652
$ws.Range("H111").Value = 59846.332
$ws.Range("J111").Value = 72636.5
$ws.Range("L111").Value = 72636.5
$ws.Range("N111").Value = -78770.5
$ws.Range("H122").Value = 4679.5654
$ws.Range("J122").Value = 5287.25
$ws.Range("L122").Value = 15861.75
$ws.Range("N122").Value = -20761.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3721.889
$ws.Range("I68").Value = 4124.25
$ws.Range("J68").Value = 3400
$ws.Range("K68").Value = 4124.25
$ws.Range("L68").Value = 3400
$ws.Range("M68").Value = -3375.25
$ws.Range("N68").Value = -4898
$ws.Range("H71").Value = 3721.889
$ws.Range("I71").Value = 4124.25
$ws.Range("J71").Value = 3400
$ws.Range("K71").Value = 20621.25
$ws.Range("L71").Value = 17000
$ws.Range("M71").Value = -16877.25
$ws.Range("N71").Value = -24488
$ws.Range("H122").Value = 4366.9165
$ws.Range("I122").Value = 3495
$ws.Range("K122").Value = 10485
$ws.Range("M122").Value = -8035
$ws.Range("H132").Value = 2691.0557
$ws.Range("I132").Value = 2486.5
$ws.Range("K132").Value = 7459.5
$ws.Range("M132").Value = -4929.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 12775
$ws.Range("I28").Value = 20000
$ws.Range("J28").Value = 5550
$ws.Range("K28").Value = 20000
$ws.Range("L28").Value = 5550
$ws.Range("M28").Value = -19652
$ws.Range("N28").Value = -6246
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40496
$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41716
$ws.Range("H112").Value = 45000
$ws.Range("J112").Value = 45000
$ws.Range("L112").Value = 45000
$ws.Range("N112").Value = -47954
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 1615.5
$ws.Range("I122").Value = 1615.5
$ws.Range("K122").Value = 4846.5
$ws.Range("M122").Value = -2396.5
$ws.Range("H132").Value = 10000.637
$ws.Range("I132").Value = 5625.25
$ws.Range("K132").Value = 16875.75
$ws.Range("M132").Value = -14345.75
$ws.Range("H136").Value = 52235.684
$ws.Range("I136").Value = 61033
$ws.Range("J136").Value = 5316.6665
$ws.Range("K136").Value = 183099
$ws.Range("L136").Value = 15949.9995
$ws.Range("M136").Value = -180549
$ws.Range("N136").Value = -21049.9995
$ws.Range("H141").Value = 82391
$ws.Range("J141").Value = 82391
$ws.Range("L141").Value = 82391
$ws.Range("N141").Value = -92751
